$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 372, shifting existing rows 372:391 down to 373:392
$ws.Rows("372:372").Insert()

# Populate the newly inserted row 372 with the new weekly price entry
$ws.Range("A372").Value = 10
$ws.Range("B372").Value = "Vega Modelo de Temuco"
$ws.Range("C372").Value = "La Araucanía"
$ws.Range("D372").Value = 44714
$ws.Range("E372").Value = 9
$ws.Range("F372").Value = 100114014
$ws.Range("G372").Value = "Betarraga"
$ws.Range("H372").Value = "Sin especificar"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 125
$ws.Range("K372").Value = 9000
$ws.Range("L372").Value = 9000
$ws.Range("M372").Value = 9000
$ws.Range("N372").Value = "`$/saco 25 kilos"
$ws.Range("O372").Value = "Provincia de Cautín"
$ws.Range("P372").Value = 360
$ws.Range("Q372").Value = 25
$ws.Range("R372").Value = "Hortaliza"
